# Update PLC data 2025-10-13 13:57:36
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LiveData")

$ws.Range("B2").Value = 7189
$ws.Range("C3").Value = 168878
$ws.Range("C4").Value = 159716
$ws.Range("C5").Value = 9162
$ws.Range("C8").Value = 65.48
